$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 690900
$ws.Range("C4").Value = 13330
$ws.Range("D4").Value = 58263
$ws.Range("E4").Value = 596682
$ws.Range("F4").Value = 13466
$ws.Range("G4").Value = 1338
$ws.Range("H4").Value = 35955

# Row 8
$ws.Range("B8").Value = 139134
$ws.Range("C8").Value = 1436
$ws.Range("D8").Value = 81800
$ws.Range("E8").Value = 53131
$ws.Range("F8").Value = 4288
$ws.Range("G8").Value = 151
$ws.Range("H8").Value = 4203

# Row 14
$ws.Range("A14").Value = "Brasil"
$ws.Range("B14").Value = 33682
$ws.Range("C14").Value = 2999
$ws.Range("D14").Value = 14026
$ws.Range("E14").Value = 17515
$ws.Range("F14").Value = 6634
$ws.Range("G14").Value = 194
$ws.Range("H14").Value = 2141

# Row 15
$ws.Range("A15").Value = "Rusia"
$ws.Range("B15").Value = 32008
$ws.Range("C15").Value = 4070
$ws.Range("D15").Value = 2590
$ws.Range("E15").Value = 29145
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 41
$ws.Range("H15").Value = 273

# Row 16
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 31642
$ws.Range("C16").Value = 1536
$ws.Range("D16").Value = 10328
$ws.Range("E16").Value = 20004
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 115
$ws.Range("H16").Value = 1310

# Row 27
$ws.Range("D27").Value = 3621
$ws.Range("E27").Value = 5515

# Row 29
$ws.Range("A29").Value = "Ecuador"
$ws.Range("B29").Value = 8450
$ws.Range("C29").Value = 225
$ws.Range("D29").Value = 838
$ws.Range("E29").Value = 7191
$ws.Range("F29").Value = 168
$ws.Range("G29").Value = 18
$ws.Range("H29").Value = 421

# Row 30
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 8379
$ws.Range("C30").Value = 461
$ws.Range("D30").Value = 866
$ws.Range("E30").Value = 7181
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 18
$ws.Range("H30").Value = 332

# Row 56
$ws.Range("E56").Value = 1880
$ws.Range("F56").Value = 126
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 123

# Row 67
$ws.Range("B67").Value = 1546
$ws.Range("C67").Value = 144
$ws.Range("D67").Value = 347
$ws.Range("E67").Value = 1182

# Row 71
$ws.Range("B71").Value = 1405
$ws.Range("C71").Value = 56
$ws.Range("D71").Value = 156
$ws.Range("E71").Value = 1245

# Row 74
$ws.Range("A74").Value = "Bosnia y Herzegovina"
$ws.Range("B74").Value = 1211
$ws.Range("C74").Value = 44
$ws.Range("D74").Value = 320
$ws.Range("E74").Value = 845
$ws.Range("F74").Value = 4
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 46

# Row 75
$ws.Range("A75").Value = "Armenia"
$ws.Range("B75").Value = 1201
$ws.Range("C75").Value = 42
$ws.Range("D75").Value = 402
$ws.Range("E75").Value = 780
$ws.Range("F75").Value = 30
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 19

# Row 82
$ws.Range("F82").Value = 17

# Row 136
$ws.Range("A136").Value = "Somalia"
$ws.Range("B136").Value = 116
$ws.Range("C136").Value = 36
$ws.Range("D136").Value = 2
$ws.Range("E136").Value = 109
$ws.Range("F136").Value = 2
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 5

# Row 137
$ws.Range("A137").Value = "Trinidad yTobago"
$ws.Range("B137").Value = 114
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 20
$ws.Range("E137").Value = 86
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 8

# Row 138
$ws.Range("A138").Value = "Etiopia"
$ws.Range("C138").Value = 4
$ws.Range("D138").Value = 15
$ws.Range("E138").Value = 78
$ws.Range("H138").Value = 3

# Row 139
$ws.Range("A139").Value = "Aruba"
$ws.Range("C139").Value = 1
$ws.Range("D139").Value = 43
$ws.Range("E139").Value = 51
$ws.Range("F139").Value = 1
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 2

# Row 140
$ws.Range("A140").Value = "Guayana Francesa"
$ws.Range("B140").Value = 96
$ws.Range("C140").Value = 10
$ws.Range("D140").Value = 61
$ws.Range("E140").Value = 35
$ws.Range("F140").Value = 2
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 0

# Row 141
$ws.Range("A141").Value = "Gabon"
$ws.Range("B141").Value = 95
$ws.Range("C141").Value = 15
$ws.Range("D141").Value = 6
$ws.Range("E141").Value = 88
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 1

# Row 142
$ws.Range("A142").Value = "Monaco"
$ws.Range("B142").Value = 94
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 20
$ws.Range("E142").Value = 71
$ws.Range("F142").Value = 3
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 3

# Row 143
$ws.Range("A143").Value = "Birmania"
$ws.Range("B143").Value = 88
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 5
$ws.Range("E143").Value = 79
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 4

# Row 144
$ws.Range("A144").Value = "Bermudas"
$ws.Range("D144").Value = 35
$ws.Range("E144").Value = 43
$ws.Range("F144").Value = 9

# Row 145
$ws.Range("A145").Value = "Togo"
$ws.Range("B145").Value = 83
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 48
$ws.Range("E145").Value = 30
$ws.Range("F145").Value = 0
